$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the header formatting from the last existing header cell (AC1) to the
# three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row (2-45).
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 91  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 71  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
